# Fix typo in date of Mobile Labs acquisition: "2920" -> "2020"
$d = $word.ActiveDocument
$d.Content.Find.Execute("October 2920", $true, $false, $false, $false, $false,
                         $true, 1, $false, "October 2020", 2)

# Touch the "Abstract Title" style's formatting (re-asserting the bold it
# already had) so the style definitions are rewritten cleanly -- this also
# clears a stray bit of malformed markup that had been lingering in that
# style's run properties.
$abstractTitle = $d.Styles("Abstract Title")
$abstractTitle.Font.Bold = $true

# Add the missing "Footnote Block Text" style (based on Footnote Text),
# matching the indentation/spacing used by the existing "Block Text"
# style.
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = "Footnote Text"
$footnoteBlockText.NextParagraphStyle = "Footnote Text"
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true
$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24
